# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets,
# mirroring a re-run of the site's data scraper (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 14017
$ws1.Range("F7").Value  = 174
$ws1.Range("F9").Value  = 57
$ws1.Range("F12").Value = 11
$ws1.Range("F13").Value = 5
$ws1.Range("F14").Value = 14301
$ws1.Range("F16").Value = 643
$ws1.Range("F17").Value = 15077
$ws1.Range("F19").Value = 8441
$ws1.Range("F30").Value = 1055
$ws1.Range("F31").Value = 1
$ws1.Range("F33").Value = 35
$ws1.Range("F36").Value = 407
$ws1.Range("F40").Value = 236
$ws1.Range("F41").Value = 399
$ws1.Range("F43").Value = 5200

# Sheet "全部类型" (All types) - row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 14017
$ws4.Range("F7").Value  = 174
$ws4.Range("F9").Value  = 57
$ws4.Range("F12").Value = 11
$ws4.Range("F13").Value = 5
$ws4.Range("F14").Value = 14301
$ws4.Range("F16").Value = 643
$ws4.Range("F17").Value = 15077
$ws4.Range("F19").Value = 8441
$ws4.Range("F31").Value = 1055
$ws4.Range("F32").Value = 1
$ws4.Range("F34").Value = 35
$ws4.Range("F39").Value = 407
$ws4.Range("F43").Value = 236
$ws4.Range("F44").Value = 399
$ws4.Range("F46").Value = 5200
